# Add a new trade record (row 16) to the sheet, mirroring the layout of
# the existing rows (Principle, Start Principle, BuyPrice, SellPrice,
# IsShortSell, Price Change %, Date, Profitable).
#
# Copy the previous row's cell formatting first (so the Date column keeps
# its existing date-time number format / style instead of creating a new
# one) and then overwrite the values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A15:H15").Copy()
$ws.Range("A16:H16").PasteSpecial()

$ws.Range("A16").Value = 9377.01
$ws.Range("B16").Value = 9732.24
$ws.Range("C16").Value = 277
$ws.Range("D16").Value = 287.11
$ws.Range("E16").Value = $true
$ws.Range("F16").Value = 3.65
$ws.Range("G16").Value = 42626.545358796298
$ws.Range("H16").Value = $false
